$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Floating Tower Properties" sheet by copying "Tower Properties" ---
# (matches header-row styles / column widths / number formats exactly), then rename + reposition.
$wsTower = $wb.Worksheets.Item("Tower Properties")
$wsMaterial = $wb.Worksheets.Item("Material Properties")
$wsTower.Copy($wsMaterial, $null) | Out-Null
$ws = $wb.Worksheets.Item("Tower Properties (2)")
$ws.Name = "Floating Tower Properties"

# Trim to 21 data rows (Tower Properties has 33 rows; Floating Tower Properties only has 21)
$ws.Range("A22:K33").EntireRow.Delete() | Out-Null

# Clear stray location labels copied over from Tower Properties (only row 2 "Tower start"
# and row 21 "Tower top" carry a label on this sheet)
$ws.Range("A3:A20").ClearContents() | Out-Null

# Overwrite with the floating-tower-specific data
$ws.Range("A2").Value = "Tower start"
$ws.Range("B2").Value = 15
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 88.528000000000006
$ws.Range("E2").Value = 21501.229098279051
$ws.Range("F2").Value = 264048.83909909509
$ws.Range("G2").Value = 264048.83909909509
$ws.Range("H2").Value = 6770483053822.9512
$ws.Range("I2").Value = 6770483053822.9512
$ws.Range("J2").Value = 5368993061681.6006
$ws.Range("K2").Value = 551313566622.53979
$ws.Range("B3").Value = 28
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = 88.528000000000006
$ws.Range("E3").Value = 21501.229098279051
$ws.Range("F3").Value = 264048.83909909509
$ws.Range("G3").Value = 264048.83909909509
$ws.Range("H3").Value = 6770483053822.9512
$ws.Range("I3").Value = 6770483053822.9512
$ws.Range("J3").Value = 5368993061681.6006
$ws.Range("K3").Value = 551313566622.53979
$ws.Range("B4").Value = 28.001000000000001
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 78.376999999999995
$ws.Range("E4").Value = 19055.301833420239
$ws.Range("F4").Value = 234486.7933995277
$ws.Range("G4").Value = 234486.7933995277
$ws.Range("H4").Value = 6012481882039.1729
$ws.Range("I4").Value = 6012481882039.1729
$ws.Range("J4").Value = 4767898132457.0635
$ws.Range("K4").Value = 488597482908.21118
$ws.Range("B5").Value = 41
$ws.Range("C5").Value = 10
$ws.Range("D5").Value = 78.376999999999995
$ws.Range("E5").Value = 19055.301833420239
$ws.Range("F5").Value = 234486.7933995277
$ws.Range("G5").Value = 234486.7933995277
$ws.Range("H5").Value = 6012481882039.1729
$ws.Range("I5").Value = 6012481882039.1729
$ws.Range("J5").Value = 4767898132457.0635
$ws.Range("K5").Value = 488597482908.21118
$ws.Range("B6").Value = 41.000999999999998
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 68.146000000000001
$ws.Range("E6").Value = 16584.988359070299
$ws.Range("F6").Value = 204506.10760939261
$ws.Range("G6").Value = 204506.10760939261
$ws.Range("H6").Value = 5243746348958.7842
$ws.Range("I6").Value = 5243746348958.7842
$ws.Range("J6").Value = 4158290854724.3159
$ws.Range("K6").Value = 425256111771.03339
$ws.Range("B7").Value = 54
$ws.Range("C7").Value = 10
$ws.Range("D7").Value = 68.146000000000001
$ws.Range("E7").Value = 16584.988359070299
$ws.Range("F7").Value = 204506.10760939261
$ws.Range("G7").Value = 204506.10760939261
$ws.Range("H7").Value = 5243746348958.7842
$ws.Range("I7").Value = 5243746348958.7842
$ws.Range("J7").Value = 4158290854724.3159
$ws.Range("K7").Value = 425256111771.03339
$ws.Range("B8").Value = 54.000999999999998
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 57.781999999999996
$ws.Range("E8").Value = 14077.331148695441
$ws.Range("F8").Value = 173944.84869346919
$ws.Range("G8").Value = 173944.84869346919
$ws.Range("H8").Value = 4460124325473.5674
$ws.Range("I8").Value = 4460124325473.5674
$ws.Range("J8").Value = 3536878590100.5391
$ws.Range("K8").Value = 360957208940.90863
$ws.Range("B9").Value = 67
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 57.781999999999996
$ws.Range("E9").Value = 14077.331148695441
$ws.Range("F9").Value = 173944.84869346919
$ws.Range("G9").Value = 173944.84869346919
$ws.Range("H9").Value = 4460124325473.5674
$ws.Range("I9").Value = 4460124325473.5674
$ws.Range("J9").Value = 3536878590100.5391
$ws.Range("K9").Value = 360957208940.90863
$ws.Range("B10").Value = 67.001000000000005
$ws.Range("C10").Value = 10
$ws.Range("D10").Value = 47.296999999999997
$ws.Range("E10").Value = 11535.04025774993
$ws.Range("F10").Value = 142830.52221336699
$ws.Range("G10").Value = 142830.52221336699
$ws.Range("H10").Value = 3662321082394.0249
$ws.Range("I10").Value = 3662321082394.0249
$ws.Range("J10").Value = 2904220618338.4619
$ws.Range("K10").Value = 295770263019.22913
$ws.Range("B11").Value = 80
$ws.Range("C11").Value = 10
$ws.Range("D11").Value = 47.296999999999997
$ws.Range("E11").Value = 11535.04025774993
$ws.Range("F11").Value = 142830.52221336699
$ws.Range("G11").Value = 142830.52221336699
$ws.Range("H11").Value = 3662321082394.0249
$ws.Range("I11").Value = 3662321082394.0249
$ws.Range("J11").Value = 2904220618338.4619
$ws.Range("K11").Value = 295770263019.22913
$ws.Range("B12").Value = 80.001000000000005
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 36.875999999999998
$ws.Range("E12").Value = 9002.9288352438489
$ws.Range("F12").Value = 111709.691065509
$ws.Range("G12").Value = 111709.691065509
$ws.Range("H12").Value = 2864351052961.77
$ws.Range("I12").Value = 2864351052961.77
$ws.Range("J12").Value = 2271430384998.6831
$ws.Range("K12").Value = 230844329108.81671
$ws.Range("B13").Value = 93
$ws.Range("C13").Value = 10
$ws.Range("D13").Value = 36.875999999999998
$ws.Range("E13").Value = 9002.9288352438489
$ws.Range("F13").Value = 111709.691065509
$ws.Range("G13").Value = 111709.691065509
$ws.Range("H13").Value = 2864351052961.77
$ws.Range("I13").Value = 2864351052961.77
$ws.Range("J13").Value = 2271430384998.6831
$ws.Range("K13").Value = 230844329108.81671
$ws.Range("B14").Value = 93.001000000000005
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 26.873000000000001
$ws.Range("E14").Value = 6567.3774435858404
$ws.Range("F14").Value = 81652.190880971466
$ws.Range("G14").Value = 81652.190880971466
$ws.Range("H14").Value = 2093645920024.9089
$ws.Range("I14").Value = 2093645920024.9089
$ws.Range("J14").Value = 1660261214579.7529
$ws.Range("K14").Value = 168394293425.27802
$ws.Range("B15").Value = 106
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 26.873000000000001
$ws.Range("E15").Value = 6567.3774435858404
$ws.Range("F15").Value = 81652.190880971466
$ws.Range("G15").Value = 81652.190880971466
$ws.Range("H15").Value = 2093645920024.9089
$ws.Range("I15").Value = 2093645920024.9089
$ws.Range("J15").Value = 1660261214579.7529
$ws.Range("K15").Value = 168394293425.27802
$ws.Range("B16").Value = 106.001
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 17.748999999999999
$ws.Range("E16").Value = 4341.5704298743331
$ws.Range("F16").Value = 54077.325967507757
$ws.Range("G16").Value = 54077.325967507757
$ws.Range("H16").Value = 1386598101730.969
$ws.Range("I16").Value = 1386598101730.969
$ws.Range("J16").Value = 1099572294672.658
$ws.Range("K16").Value = 111322318714.7265
$ws.Range("B17").Value = 119
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 17.748999999999999
$ws.Range("E17").Value = 4341.5704298743331
$ws.Range("F17").Value = 54077.325967507757
$ws.Range("G17").Value = 54077.325967507757
$ws.Range("H17").Value = 1386598101730.969
$ws.Range("I17").Value = 1386598101730.969
$ws.Range("J17").Value = 1099572294672.658
$ws.Range("K17").Value = 111322318714.7265
$ws.Range("B18").Value = 119.001
$ws.Range("C18").Value = 10
$ws.Range("D18").Value = 9.9109999999999996
$ws.Range("E18").Value = 2426.2263151018619
$ws.Range("F18").Value = 30267.772697042441
$ws.Range("G18").Value = 30267.772697042441
$ws.Range("H18").Value = 776096735821.6012
$ws.Range("I18").Value = 776096735821.6012
$ws.Range("J18").Value = 615444711506.52966
$ws.Range("K18").Value = 62210931156.457993
$ws.Range("B19").Value = 132
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 9.9109999999999996
$ws.Range("E19").Value = 2426.2263151018619
$ws.Range("F19").Value = 30267.772697042441
$ws.Range("G19").Value = 30267.772697042441
$ws.Range("H19").Value = 776096735821.6012
$ws.Range("I19").Value = 776096735821.6012
$ws.Range("J19").Value = 615444711506.52966
$ws.Range("K19").Value = 62210931156.457993
$ws.Range("B20").Value = 132.001
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 7.9359999999999999
$ws.Range("E20").Value = 1943.127694419381
$ws.Range("F20").Value = 24250.575121376562
$ws.Range("G20").Value = 24250.575121376562
$ws.Range("H20").Value = 621809618496.83496
$ws.Range("I20").Value = 621809618496.83496
$ws.Range("J20").Value = 493095027467.99011
$ws.Range("K20").Value = 49823787036.394386
$ws.Range("A21").Value = "Tower top"
$ws.Range("B21").Value = 144.386
$ws.Range("C21").Value = 6.5
$ws.Range("D21").Value = 7.9359999999999999
$ws.Range("E21").Value = 1262.4928495596951
$ws.Range("F21").Value = 6651.2791319294456
$ws.Range("G21").Value = 6651.2791319294456
$ws.Range("H21").Value = 170545618767.42169
$ws.Range("I21").Value = 170545618767.42169
$ws.Range("J21").Value = 135242675682.5654
$ws.Range("K21").Value = 32371611527.171661
$ws.Range("F25").Select() | Out-Null

# --- 2. Add the new Overview rows (floating substructure summary) ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A37").Value = "Floating tower mass [t]"
$ov.Range("B37").Value = 1483.073634074095
$ov.Range("A38").Value = "Floating tower base diameter [m]"
$ov.Range("B38").Value = 10
$ov.Range("A39").Value = "Floating transition piece height [m]"
$ov.Range("B39").Value = 15
$ov.Range("A40").Value = "Volturn-S hull mass [t]"
$ov.Range("B40").Value = 4014.227616744387
$ov.Range("A41").Value = "Volturn-S fixed ballast mass [t]"
$ov.Range("B41").Value = 2539.9999950000001
$ov.Range("A42").Value = "Volturn-S fluid ballast mass [t]"
$ov.Range("B42").Value = 8444.6854659379842
$ov.Range("A43").Value = "Volturn-S displacement [m^3]"
$ov.Range("B43").Value = 17755.490038649681
$ov.Range("A44").Value = "Volturn-S freeboard [m]"
$ov.Range("B44").Value = 15
$ov.Range("A45").Value = "Volturn-S draft [m]"
$ov.Range("B45").Value = 20
# --- 3. Make Overview the active sheet/selection, matching the saved view state ---
$ov.Activate() | Out-Null
$ov.Range("B38").Select() | Out-Null

# --- 4. Best-effort restore of the window chrome size/position ---
$win = $excel.ActiveWindow
$win.Left = 3580
$win.Top = 500
$win.Width = 25940
$win.Height = 21300
